# Regenerate the scraped "想去人数" (interest count) figures for a handful
# of events across the 展览 / 演出 / 全部类型 sheets, matching the refreshed
# gh-pages data output.

$wb = $excel.ActiveWorkbook

$ws_exhibit = $wb.Worksheets.Item("展览")
$ws_exhibit.Range("F5").Value = 1044
$ws_exhibit.Range("F13").Value = 525
$ws_exhibit.Range("F14").Value = 1692
$ws_exhibit.Range("F22").Value = 81
$ws_exhibit.Range("F25").Value = 3600
$ws_exhibit.Range("F26").Value = 713
$ws_exhibit.Range("F28").Value = 1591
$ws_exhibit.Range("F29").Value = 52

$ws_show = $wb.Worksheets.Item("演出")
$ws_show.Range("F7").Value = 9
$ws_show.Range("F8").Value = 35

$ws_all = $wb.Worksheets.Item("全部类型")
$ws_all.Range("F12").Value = 9
$ws_all.Range("F13").Value = 35
$ws_all.Range("F16").Value = 1044
$ws_all.Range("F24").Value = 525
$ws_all.Range("F25").Value = 1692
$ws_all.Range("F35").Value = 82
$ws_all.Range("F38").Value = 3600
$ws_all.Range("F39").Value = 713
$ws_all.Range("F41").Value = 1591
$ws_all.Range("F44").Value = 52
